# Generate Report for Handback
# - Marks both localized-file rows as handed back (status text change)
# - Records the handback target file, handback xlf file, and handback
#   datetime for zh-cn and de-de
# - Widens the Status / Latest Target File / Latest Handback File columns
#   so the longer text fits

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$repoBlobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/76b03d6bfc3a38bde41be184a7a4afbe377f4713/e2e/"
$file1 = "91f76a80-f8c5-4cb4-8657-33d9a94cbc75.md"
$file2 = "d3752f80-a8a2-4ba7-a61b-92a86e310393.md"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-language status cells + widen columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# Per-language detail sheets (zh-cn / de-de)
# ---------------------------------------------------------------------
$langSheets = @{
    "zh-cn" = @{
        XliffSuffix   = "8f0240dcf7bcc81aa1992b9d90c8ae0160c8530a.zh-cn.xlf"
        XliffSuffix2  = "f771250bc0e46d5bb54379b8e3a45dd991f2c4a7.zh-cn.xlf"
        HandbackTime  = "2016-09-01 21:11:07"
    }
    "de-de" = @{
        XliffSuffix   = "8f0240dcf7bcc81aa1992b9d90c8ae0160c8530a.de-de.xlf"
        XliffSuffix2  = "f771250bc0e46d5bb54379b8e3a45dd991f2c4a7.de-de.xlf"
        HandbackTime  = "2016-09-01 21:11:18"
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $langSheets[$sheetName]

    # Status column (C) for both rows
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Row 2 -> file1 (91f76a80...)
    $handoff1 = $file1 -replace "\.md$", "." + $info.XliffSuffix
    $ws.Hyperlinks.Add($ws.Range("I2"), $repoBlobBase + $file1, "", "", $file1)
    $ws.Range("J2").Value = $handoff1
    $ws.Range("K2").Value = $info.HandbackTime

    # Row 3 -> file2 (d3752f80...)
    $handoff2 = $file2 -replace "\.md$", "." + $info.XliffSuffix2
    $ws.Hyperlinks.Add($ws.Range("I3"), $repoBlobBase + $file2, "", "", $file2)
    $ws.Range("J3").Value = $handoff2
    $ws.Range("K3").Value = $info.HandbackTime

    # Widen columns: Status (C), Latest Target File (I), Latest Handback File (J)
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}
